# Applies the crypto-price/volume refresh described in the commit:
# "Updated cryptos list on Sun Jun 16 23:42:32 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values (column D) must stay literal text (e.g. "0.160", "8.40")
# instead of being auto-coerced to numbers, which would drop trailing/
# grouping-like characters. Force text format on every Price cell we touch.
$priceCells = @("D2","D3","D5","D6","D7","D13","D14","D15","D16","D17","D19","D20","D22","D23","D24","D26","D27","D28","D29","D31","D32","D33","D34","D35","D37","D39","D40","D41","D42","D47","D49","D50")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

# Standard per-cell updates (price / volume columns) based on the diff
$ws.Range("D2").Value = "66.663.05"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.618.06"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "610.51"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "150.47"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").Value = "3.616.46"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "4.235.22"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "0.0000209"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "29.97"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "3.618.05"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "66.749.59"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "11.70"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").Value = "427.97"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "0.619"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "78.86"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "8.40"
$ws.Range("E27").Value = "  +5.56%  "
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").Value = "  +5.71%  "
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D33").Value = "1.47"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").Value = "25.46"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "7.89"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "177.18"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "0.0862"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "0.901"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E45").Value = "  +7.77%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "25.06"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").Value = "23.94"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "7.21"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +2.28%  "

# Rows 31 and 32 swap their Coin/Link/Price data (Kaspa <-> RenzoRestakedETH),
# each also getting an updated Price/Volume figure
$ws.Range("B31").Value = "RenzoRestakedETH"
$ws.Range("C31").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D31").Value = "3.619.57"
$ws.Range("E31").Value = "  +1.39%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "0.160"
$ws.Range("E32").Value = "  +3.77%  "
